$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers section (rows 3-6) ---
$ws.Range("A3").Value = 'Intel(R) Wireless-AC 9260 160MHz - 23.40.0.4'
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 365
$ws.Range("D3").Value = 97.5
$ws.Range("A4").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.70.25.2'
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 1627
$ws.Range("D4").Value = 97.8
$ws.Range("A5").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 23.120.0.3'
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 98
$ws.Range("D5").Value = 98.6
$ws.Range("A6").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.250.0.4'
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 1021
$ws.Range("D6").Value = 98.7

# --- Totals row 7 ---
$ws.Range("C7").Value = 3111

# --- Good Drivers section (rows 15-49): A, B, D columns ---
$ws.Range("A15").Value = 'Intel(R) Wireless-AC 9260 160MHz - 23.120.0.3'
$ws.Range("B15").Value = 17963
$ws.Range("D15").Value = 100
$ws.Range("A16").Value = 'Intel(R) Wireless-AC 9260 160MHz - 23.100.0.4'
$ws.Range("B16").Value = 108115
$ws.Range("D16").Value = 99.9
$ws.Range("A17").Value = 'Intel(R) Wireless-AC 9260 160MHz - 23.80.0.7'
$ws.Range("B17").Value = 136392
$ws.Range("D17").Value = 100
$ws.Range("A18").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 23.70.2.3'
$ws.Range("B18").Value = 18721
$ws.Range("D18").Value = 99.9
$ws.Range("A19").Value = 'Intel(R) Wireless-AC 9260 160MHz - 23.60.1.2'
$ws.Range("B19").Value = 45212
$ws.Range("D19").Value = 99.9
$ws.Range("A20").Value = 'Intel(R) Wireless-AC 9260 160MHz - 23.50.0.6'
$ws.Range("B20").Value = 787603
$ws.Range("D20").Value = 99.9
$ws.Range("A21").Value = 'Intel(R) Wireless-AC 9260 160MHz - 23.20.1.1'
$ws.Range("B21").Value = 32457
$ws.Range("D21").Value = 100
$ws.Range("A22").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.250.10.1'
$ws.Range("B22").Value = 69578
$ws.Range("D22").Value = 99.9
$ws.Range("A23").Value = 'Intel(R) Wireless-AC 9260 160MHz - 22.240.0.6'
$ws.Range("B23").Value = 100154
$ws.Range("D23").Value = 100
$ws.Range("A24").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.230.0.8'
$ws.Range("B24").Value = 329845
$ws.Range("D24").Value = 99.9
$ws.Range("A25").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.200.0.6'
$ws.Range("B25").Value = 143808
$ws.Range("D25").Value = 99.9
$ws.Range("A26").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.190.0.4'
$ws.Range("B26").Value = 287148
$ws.Range("D26").Value = 99.9
$ws.Range("A27").Value = 'Intel(R) Wireless-AC 9260 160MHz - 22.180.0.4'
$ws.Range("B27").Value = 89463
$ws.Range("D27").Value = 100
$ws.Range("A28").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.160.0.4'
$ws.Range("B28").Value = 96526
$ws.Range("D28").Value = 99.9
$ws.Range("A29").Value = 'Intel(R) Wireless-AC 9260 160MHz - 22.160.0.4'
$ws.Range("B29").Value = 56294
$ws.Range("D29").Value = 100
$ws.Range("A30").Value = 'Intel(R) Wireless-AC 9260 160MHz - 22.130.0.5'
$ws.Range("B30").Value = 34662
$ws.Range("D30").Value = 100
$ws.Range("A31").Value = 'Intel(R) Wireless-AC 9260 160MHz - 22.120.0.3'
$ws.Range("B31").Value = 116879
$ws.Range("D31").Value = 100
$ws.Range("A32").Value = 'Intel(R) Wireless-AC 9260 160MHz - 22.40.0.7'
$ws.Range("B32").Value = 103051
$ws.Range("D32").Value = 100
$ws.Range("A33").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 22.30.0.11'
$ws.Range("B33").Value = 170510
$ws.Range("D33").Value = 99.9
$ws.Range("A34").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.30.0.11'
$ws.Range("B34").Value = 67111
$ws.Range("D34").Value = 100
$ws.Range("A35").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.10.0.7'
$ws.Range("B35").Value = 66577
$ws.Range("D35").Value = 100
$ws.Range("A36").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.0.1.1'
$ws.Range("B36").Value = 15734
$ws.Range("D36").Value = 99.9
$ws.Range("A37").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 22.0.1.1'
$ws.Range("B37").Value = 52096
$ws.Range("D37").Value = 100
$ws.Range("A38").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 21.60.2.1'
$ws.Range("B38").Value = 26241
$ws.Range("D38").Value = 100
$ws.Range("A39").Value = 'Intel(R) Wireless-AC 9260 160MHz - 21.60.0.5'
$ws.Range("B39").Value = 54452
$ws.Range("D39").Value = 100
$ws.Range("A40").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.70.11.3'
$ws.Range("B40").Value = 161874
$ws.Range("D40").Value = 100
$ws.Range("A41").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 21.40.2.2'
$ws.Range("B41").Value = 88435
$ws.Range("D41").Value = 99.9
$ws.Range("A42").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.70.12.5'
$ws.Range("B42").Value = 143342
$ws.Range("D42").Value = 99.9
$ws.Range("A43").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 21.30.4.1'
$ws.Range("B43").Value = 13016
$ws.Range("D43").Value = 100
$ws.Range("A44").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.70.10.2'
$ws.Range("B44").Value = 20227
$ws.Range("D44").Value = 100
$ws.Range("A45").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.70.9.1'
$ws.Range("B45").Value = 34065
$ws.Range("D45").Value = 100
$ws.Range("A46").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 21.10.1.2'
$ws.Range("B46").Value = 46270
$ws.Range("D46").Value = 100
$ws.Range("A47").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.70.8.1'
$ws.Range("B47").Value = 48540
$ws.Range("D47").Value = 100
$ws.Range("A48").Value = 'Intel(R) Wireless-AC 9260 160MHz - 20.120.0.100'
$ws.Range("B48").Value = 55058
$ws.Range("D48").Value = 100
$ws.Range("A49").Value = 'Intel(R) Dual Band Wireless-AC 8265 - 20.70.5.2'
$ws.Range("B49").Value = 184564
$ws.Range("D49").Value = 99.9

# --- Good Drivers section E column (Driver Vintage dates, stored as text) ---
# Set raw text values with apostrophe prefix to force text interpretation
$ws.Range("E15").Value = "'2025-02-05"
$ws.Range("E16").Value = "'2024-11-10"
$ws.Range("E17").Value = "'2024-09-01"
$ws.Range("E18").Value = "'2024-07-23"
$ws.Range("E19").Value = "'2024-06-02"
$ws.Range("E20").Value = "'2024-04-13"
$ws.Range("E21").Value = "'2023-12-19"
$ws.Range("E22").Value = "'2023-08-14"
$ws.Range("E23").Value = "'2023-06-17"
$ws.Range("E24").Value = "'2023-05-08"
$ws.Range("E25").Value = "'2023-01-16"
$ws.Range("E26").Value = "'2022-11-22"
$ws.Range("E27").Value = "'2022-10-17"
$ws.Range("E28").Value = "'2022-08-13"
$ws.Range("E29").Value = "'2022-08-13"
$ws.Range("E30").Value = "'2022-03-14"
$ws.Range("E31").Value = "'2022-01-30"
$ws.Range("E32").Value = "'2021-03-02"
$ws.Range("E33").Value = "'2021-01-19"
$ws.Range("E34").Value = "'2021-01-19"
$ws.Range("E35").Value = "'2020-10-19"
$ws.Range("E36").Value = "'2020-09-28"
$ws.Range("E37").Value = "'2020-09-28"
$ws.Range("E38").Value = "'2019-12-14"
$ws.Range("E39").Value = "'2019-11-10"
$ws.Range("E40").Value = "'2019-09-05"
$ws.Range("E41").Value = "'2019-08-31"
$ws.Range("E42").Value = "'2019-08-25"
$ws.Range("E43").Value = "'2019-07-29"
$ws.Range("E44").Value = "'2019-05-11"
$ws.Range("E45").Value = "'2019-04-28"
$ws.Range("E46").Value = "'2019-04-23"
$ws.Range("E47").Value = "'2019-03-16"
$ws.Range("E48").Value = "'2019-01-27"
$ws.Range("E49").Value = "'2018-11-25"

# Copy the number format of D15 (style matches original s=4) onto E15:E49
$ws.Range("D15").Copy()
$ws.Range("E15:E49").PasteSpecial(-4122)
$excel.CutCopyMode = $false
